$data = @(
    @("Driffin","Daniel","dnc.driffin.txt","Wednesday","speech","Cision"),
    @("Tanden","Neera","dnc.tanden.txt","Wednesday","speech","Cision"),
    @("Grisham","Michelle","dnc.grisham.txt","Wednesday","speech","Cision"),
    @("Norton","Eleanor","dnc.norton.txt","Wednesday","speech","Cision"),
    @("Schiff","Adam","dnc.schiff.txt","Wednesday","speech","Cision"),
    @("Waters","Maxine","dnc.waters.txt","Wednesday","speech","Cision"),
    @("Hogue","Ilyse","dnc.hogue.txt","Wednesday","speech","Cision"),
    @("Gillum","Andrew","dnc.gillum.txt","Wednesday","speech","Cision"),
    @("Asian","Caucus","dnc.asian.txt","Wednesday","speech","Cision"),
    @("Bell","Brooks","dnc.bellbrooks.txt","Wednesday","speech","Cision"),
    @("deBlasio","Bill","dnc.deblasio.txt","Wednesday","speech","Cision"),
    @("Grubbe","Jeff","dnc.grubbe.txt","Wednesday","speech","Cision"),
    @("Jackson","Jesse","dnc.jacksonjesse.txt","Wednesday","speech","Cision"),
    @("Jones","Star","dnc.jonesstar.txt","Wednesday","speech","Cision"),
    @("Weaver","Karen","dnc.weaver.txt","Wednesday","speech","Cision"),
    @("Black","Caucus","dnc.blackcaucus.txt","Wednesday","speech","Cision"),
    @("Schriock","tephanie","dnc.schriock.txt","Wednesday","speech","Cision"),
    @("Reid","Harry","dnc.reidharry.txt","Wednesday","speech","Cision"),
    @("Dorff","jamie","dnc.dorff.txt","Wednesday","speech","Cision"),
    @("Belkofer","Sharon","dnc.belkofer.txt","Wednesday","speech","Cision"),
    @("Salguero","Gabriel","dnc.salguero.txt","Wednesday","speech","Cision"),
    @("Lee","Shelia","dnc.leeshelia.txt","Wednesday","benediction","Cision"),
    @("Duggan","Mike","dnc.duggan.txt","Wednesday","speech","Cision"),
    @("Omalley","Martin","dnc.omalley.txt","Wednesday","speech","Cision"),
    @("Weaver","Sigourney","dnc.weaversigourney.txt","Wednesday","speech","Cision"),
    @("Brown","Jerry","dnc.brownjerry.txt","Wednesday","speech","Cision"),
    @("Daniels","Lee","dnc.danielslee.txt","Wednesday","speech","Cision"),
    @("Murphy","Chris","dnc.murphychris.txt","Wednesday","speech","Cision"),
    @("Ramsey","Charles","dnc.ramsey.txt","Wednesday","speech","Cision"),
    @("Bassett","Angela","dnc.bassett.txt","Wednesday","speech","Cision"),
    @("Sanders","Felicia","dnc.sandersfelicia.txt","Wednesday","speech","Cision"),
    @("Sheppard","polly","dnc.sheppard.txt","Wednesday","speech","Cision"),
    @("Kelly","Mark","dnc.kellymark.txt","Wednesday","speech","Cision"),
    @("Hutson","John","dnc.hutson.txt","Wednesday","speech","Cision"),
    @("Kavanaugh","Kristen","dnc.kavanaugh.txt","Wednesday","speech","Cision"),
    @("Panetta","Leon","dnc.panetta.txt","Wednesday","speech","Cision"),
    @("Lujan","Ben","dnc.lujan.txt","Wednesday","speech","Cision"),
    @("Smegielski","Erica","dnc.smegielski.txt","Wednesday","speech","Cision"),
    @("Reed","Kasim","dnc.reedkasim.txt","Wednesday","speech","Cision"),
    @("Scott","Robert","dnc.scottbobby.txt","Wednesday","speech","Cision"),
    @("America","Our","dnc.ouramerica.txt","Wednesday","speech","Cision"),
)

$wb = $excel.ActiveWorkbook

# --- Populate the new Wednesday rows (rows 9-49) ---
$wsWed = $wb.Worksheets.Item("Wednesday")
$startRow = 9
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    for ($c = 0; $c -lt $vals.Length; $c++) {
        $wsWed.Cells.Item($row, $c + 1).Value = $vals[$c]
    }
}

# --- Update view state: Wednesday becomes the active/selected sheet ---
$wsTue = $wb.Worksheets.Item("Tuesday")
$wsTue.Activate()
$tueWin = $excel.ActiveWindow
$tueWin.SplitRow = 1
$tueWin.FreezePanes = $true
$tueWin.ScrollRow = 27
$tueWin.RangeSelection.Select() | Out-Null

$wsWed.Activate()
$wedWin = $excel.ActiveWindow
$wedWin.Zoom = 125
$wedWin.SplitRow = 1
$wedWin.FreezePanes = $true
$wedWin.ScrollRow = 39
$wsWed.Range("E50").Select()
